# ---------------------------------------------------------------------------
# "final optimization from template"
#
# - adds a "stage" column to the front of "configuration" and wires cell
#   D2 up to a VLOOKUP() against a new "urls" sheet, driven by a dropdown
#   (list data validation) on A2 sourced from urls!$A$1:$A$3
# - adds the new "urls" worksheet (qa/dev/prod -> environment URL table)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "configuration": insert a new leading "stage" column (A) - this shifts
#    every other column right by one and keeps data/format/validations in
#    lock-step.
# ---------------------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("configuration")
$wsConfig.Columns("A:A").Insert()

$wsConfig.Range("A1").Value = "stage"
$wsConfig.Range("A2").Value = "qa"

# ---------------------------------------------------------------------------
# 2. Build the new "urls" lookup sheet (inserted right after "login")
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("login")
$wsUrls = $wb.Worksheets.Add($null, $wsLogin)
$wsUrls.Name = "urls"

$wsUrls.Range("A1").Value = "qa"
$wsUrls.Range("A2").Value = "dev"
$wsUrls.Range("A3").Value = "prod"
$wsUrls.Range("B1").Value = "https://rahulshettyacademy.com/angularpractice/"
$wsUrls.Range("B2").Value = "https://dev.rahulshettyacademy.com/angularpractice/"
$wsUrls.Range("B3").Value = "https://prod.rahulshettyacademy.com/angularpractice/"

$wsUrls.Hyperlinks.Add($wsUrls.Range("B1"), "https://rahulshettyacademy.com/angularpractice/")
$wsUrls.Hyperlinks.Add($wsUrls.Range("B2"), "https://dev.rahulshettyacademy.com/angularpractice/")
$wsUrls.Hyperlinks.Add($wsUrls.Range("B3"), "https://prod.rahulshettyacademy.com/angularpractice/")
$wsUrls.Range("B1:B3").Style = "Hyperlink"

$wsUrls.Columns("B:B").ColumnWidth = 46
[void]$wsUrls.Range("B4").Select()

# ---------------------------------------------------------------------------
# 3. Back on "configuration": re-anchor the hyperlink that got left behind
#    on its old (pre-insert) address, wire up the VLOOKUP formula and the
#    urls-backed dropdown on A2.
# ---------------------------------------------------------------------------

# the hyperlink on the (now) api_url cell kept pointing at its old B2
# address across the column insert - drop it and re-anchor on C2.
$wsConfig.Hyperlinks.Delete()
$wsConfig.Hyperlinks.Add($wsConfig.Range("C2"), "https://api.rahulshettyacademy.com/angularpractice/")
$wsConfig.Range("C2").Style = "Hyperlink"

# frontend_url (D2) becomes a lookup against the new "urls" sheet, keyed by
# the "stage" dropdown in A2.
$wsConfig.Range("D2").Formula = "=VLOOKUP(A2,urls!A1:B3,2,FALSE)"

# dropdown list on A2, sourced from the "urls" sheet's first column.
$wsConfig.Range("A2").Validation.Delete()
$wsConfig.Range("A2").Validation.Add(3, 1, 1, "=urls!`$A`$1:`$A`$3")
